$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Venom Energy Special Edition / Naruto Uzomaki" record was mis-filed
# under the "RS" sales area (row 37). Remove it from there...
$ws.Rows("37:37").Delete()

# ...and re-insert it in its correct spot, under the "RN1" sales area,
# pushing the remaining RN1..RW rows back down by one.
$ws.Rows("19:19").Insert()
$ws.Range("A19").Value = "RN1"
$ws.Range("B19").Value = "Venom Energy Special Edition"
$ws.Range("C19").Value = "Naruto Uzomaki"

# Restore the user's selection to the cell they left active.
$ws.Range("B14").Select()
